$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.300.36"
$ws.Range("E2").Value = "  +2.64%  "
$ws.Range("D3").Value = "1.899.16"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("E4").Value = "  -1.17%  "
$ws.Range("D5").Value = "'315.22"
$ws.Range("E5").Value = "  -0.62%  "
$ws.Range("E6").Value = "  -1.09%  "
$ws.Range("D7").Value = "'0.5150"
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("D8").Value = "'0.3931"
$ws.Range("E8").Value = "  -1.07%  "
$ws.Range("D9").Value = "'0.08444"
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("D10").Value = "'42.49"
$ws.Range("E10").Value = "  +1.46%  "
$ws.Range("D11").Value = "'1.117"
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("D12").Value = "'6.252"
$ws.Range("E12").Value = "  -0.63%  "
$ws.Range("D13").Value = "1.895.39"
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("E14").Value = "  +0.96%  "
$ws.Range("D15").Value = "'7.325"
$ws.Range("E15").Value = "  +0.31%  "
$ws.Range("D16").Value = "'1.004"
$ws.Range("E16").Value = "  -1.20%  "
$ws.Range("D17").Value = "'93.24"
$ws.Range("E17").Value = "  +1.92%  "
$ws.Range("D18").Value = "'0.00001106"
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("D19").Value = "'0.06739"
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("E21").Value = "  -1.11%  "
$ws.Range("D22").Value = "'6.028"
$ws.Range("E22").Value = "  +0.98%  "
$ws.Range("D23").Value = "29.305.23"
$ws.Range("E23").Value = "  +2.54%  "
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("E25").Value = "  -2.59%  "
$ws.Range("D26").Value = "2.112.58"
$ws.Range("E26").Value = "  +0.42%  "
$ws.Range("D27").Value = "'159.18"
$ws.Range("E27").Value = "  -1.63%  "
$ws.Range("E28").Value = "  +0.80%  "
$ws.Range("D29").Value = "'2.437"
$ws.Range("E29").Value = "  +1.68%  "
$ws.Range("D30").Value = "'128.36"
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("D31").Value = "'1.062"
$ws.Range("E31").Value = "  +0.80%  "
$ws.Range("E32").Value = "  -1.07%  "
$ws.Range("D33").Value = "'6.177"
$ws.Range("E33").Value = "  +6.44%  "
$ws.Range("D34").Value = "'3.657"
$ws.Range("E34").Value = "  +0.92%  "
$ws.Range("D35").Value = "'0.02473"
$ws.Range("E35").Value = "  +1.00%  "
$ws.Range("D36").Value = "'0.06562"
$ws.Range("E36").Value = "  +0.59%  "
$ws.Range("D37").Value = "'9.068"
$ws.Range("E37").Value = "  +1.48%  "
$ws.Range("D38").Value = "'0.2195"
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("D39").Value = "'1.235"
$ws.Range("E39").Value = "  +3.61%  "
$ws.Range("D40").Value = "'5.122"
$ws.Range("E40").Value = "  +1.47%  "
$ws.Range("D41").Value = "'0.6499"
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").Value = "'1.237"
$ws.Range("E42").Value = "  -2.71%  "
$ws.Range("E43").Value = "  +0.25%  "
$ws.Range("D44").Value = "'0.6069"
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("D45").Value = "'13.18"
$ws.Range("E45").Value = "  +0.70%  "
$ws.Range("D46").Value = "'3.674"
$ws.Range("E46").Value = "  -1.24%  "
$ws.Range("E47").Value = "  +1.27%  "
$ws.Range("E48").Value = "  +1.15%  "
$ws.Range("D49").Value = "'123.42"
$ws.Range("E49").Value = "  +0.44%  "
$ws.Range("D50").Value = "'1.175"
$ws.Range("E50").Value = "  -2.53%  "
$ws.Range("D51").Value = "'77.77"
$ws.Range("E51").Value = "  +0.67%  "
